# daily auto push: 2026-02-07 02:49 UTC
#
# A new sample row (2026/02/07, 土, 8, 201) was recorded for the running log
# on sheet "Sheet1". It belongs right after the existing 2026/02/07 row
# (row 794), so every subsequent row shifts down by one and the new
# row's data is written into the freshly opened row 794.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one blank row at row 794, shifting rows 794:835 down to 795:836.
$ws.Rows.Item(794).Insert()

# Column A holds dates formatted as plain text (e.g. "2026/12/29"), not
# real date serials. Force text formatting before assigning so the
# "2026/02/07"-looking value isn't auto-coerced into a date number.
$ws.Cells.Item(794, 1).NumberFormat = "@"
$ws.Cells.Item(794, 1).Value = "2026/02/07"
$ws.Cells.Item(794, 2).Value = "土"
$ws.Cells.Item(794, 3).Value = 8
$ws.Cells.Item(794, 4).Value = 201
